$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 705.93335
$ws.Range("J32").Value = 924.8
$ws.Range("L32").Value = 924.8
$ws.Range("N32").Value = -1576.8
$ws.Range("H41").Value = 15625441
$ws.Range("I41").Value = 15625441
$ws.Range("K41").Value = 15625441
$ws.Range("M41").Value = -15625001
$ws.Range("H70").Value = 2176.4546
$ws.Range("J70").Value = 1934
$ws.Range("L70").Value = 5802
$ws.Range("N70").Value = -6342
$ws.Range("H73").Value = 2176.4546
$ws.Range("J73").Value = 1934
$ws.Range("L73").Value = 5802
$ws.Range("N73").Value = -7674
$ws.Range("H80").Value = 46668
$ws.Range("I80").Value = 33649.668
$ws.Range("K80").Value = 100949.004
$ws.Range("M80").Value = -99951.00399999999
$ws.Range("H83").Value = 46668
$ws.Range("I83").Value = 33649.668
$ws.Range("K83").Value = 302847.012
$ws.Range("M83").Value = -297855.012
$ws.Range("H98").Value = 8774.714
$ws.Range("I98").Value = 8570.666999999999
$ws.Range("J98").Value = 9999
$ws.Range("K98").Value = 8570.666999999999
$ws.Range("L98").Value = 9999
$ws.Range("M98").Value = -7072.666999999999
$ws.Range("N98").Value = -12995
$ws.Range("H107").Value = 62501524
$ws.Range("I107").Value = 25001828
$ws.Range("K107").Value = 25001828
$ws.Range("M107").Value = -24999908
$ws.Range("H122").Value = 8774.714
$ws.Range("I122").Value = 8570.666999999999
$ws.Range("J122").Value = 9999
$ws.Range("K122").Value = 25712.001
$ws.Range("L122").Value = 29997
$ws.Range("M122").Value = -23262.001
$ws.Range("N122").Value = -34897
$ws.Range("H132").Value = 2081.087
$ws.Range("I132").Value = 2081.087
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 6243.261
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -3713.261
$ws.Range("N132").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 4850.773
$ws.Range("I45").Value = 1181.909
$ws.Range("J45").Value = 8519.637000000001
$ws.Range("K45").Value = 1181.909
$ws.Range("L45").Value = 8519.637000000001
$ws.Range("M45").Value = -804.9090000000001
$ws.Range("N45").Value = -9273.637000000001
$ws.Range("H61").Value = 5204.241
$ws.Range("I61").Value = 1597.65
$ws.Range("K61").Value = 1597.65
$ws.Range("M61").Value = -1385.65
$ws.Range("H74").Value = 94452.17999999999
$ws.Range("I74").Value = 144996.28
$ws.Range("J74").Value = 6000
$ws.Range("K74").Value = 144996.28
$ws.Range("L74").Value = 6000
$ws.Range("M74").Value = -144122.28
$ws.Range("N74").Value = -7748
$ws.Range("H77").Value = 94452.17999999999
$ws.Range("I77").Value = 144996.28
$ws.Range("J77").Value = 6000
$ws.Range("K77").Value = 724981.4
$ws.Range("L77").Value = 30000
$ws.Range("M77").Value = -720613.4
$ws.Range("N77").Value = -38736
$ws.Range("H102").Value = 4539.8335
$ws.Range("I102").Value = 4320.5557
$ws.Range("K102").Value = 4320.5557
$ws.Range("M102").Value = -2698.5557
$ws.Range("H110").Value = 42918788
$ws.Range("I110").Value = 2502174
$ws.Range("J110").Value = 83335400
$ws.Range("K110").Value = 2502174
$ws.Range("L110").Value = 83335400
$ws.Range("M110").Value = -2500129
$ws.Range("N110").Value = -83339490
$ws.Range("H132").Value = 6024.44
$ws.Range("I132").Value = 3744.1765
$ws.Range("K132").Value = 11232.5295
$ws.Range("M132").Value = -8702.529500000001
$ws.Range("H136").Value = 5204.241
$ws.Range("I136").Value = 1597.65
$ws.Range("K136").Value = 4792.950000000001
$ws.Range("M136").Value = -2242.950000000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 632.9375
$ws.Range("I94").Value = 407.9091
$ws.Range("K94").Value = 407.9091
$ws.Range("M94").Value = 43.09089999999998
$ws.Range("H99").Value = 8266341
$ws.Range("I99").Value = 883
$ws.Range("J99").Value = 22730892
$ws.Range("K99").Value = 883
$ws.Range("L99").Value = 22730892
$ws.Range("M99").Value = 615
$ws.Range("N99").Value = -22733888

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 4415.1816
$ws.Range("I16").Value = 2424.2
$ws.Range("J16").Value = 6074.3335
$ws.Range("K16").Value = 2424.2
$ws.Range("L16").Value = 6074.3335
$ws.Range("M16").Value = -2137.2
$ws.Range("N16").Value = -6648.3335
$ws.Range("H103").Value = 26673
$ws.Range("I103").Value = 5500
$ws.Range("J103").Value = 37259.5
$ws.Range("K103").Value = 5500
$ws.Range("L103").Value = 37259.5
$ws.Range("M103").Value = -4328
$ws.Range("N103").Value = -39603.5
$ws.Range("H113").Value = 4415.1816
$ws.Range("I113").Value = 2424.2
$ws.Range("J113").Value = 6074.3335
$ws.Range("K113").Value = 2424.2
$ws.Range("L113").Value = 6074.3335
$ws.Range("M113").Value = -254.1999999999998
$ws.Range("N113").Value = -10414.3335
$ws.Range("H134").Value = 7957.2
$ws.Range("I134").Value = 7050.44
$ws.Range("J134").Value = 9090.65
$ws.Range("K134").Value = 21151.32
$ws.Range("L134").Value = 27271.95
$ws.Range("M134").Value = -18616.32
$ws.Range("N134").Value = -32341.95

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 5228.2144
$ws.Range("J34").Value = 6634.909
$ws.Range("L34").Value = 19904.727
$ws.Range("N34").Value = -20072.727
$ws.Range("H39").Value = 11171.777
$ws.Range("I39").Value = 833.3333
$ws.Range("J39").Value = 16341
$ws.Range("K39").Value = 2499.9999
$ws.Range("L39").Value = 49023
$ws.Range("M39").Value = -2205.9999
$ws.Range("N39").Value = -49611
$ws.Range("H55").Value = 48156204
$ws.Range("J55").Value = 20013600
$ws.Range("L55").Value = 60040800
$ws.Range("N55").Value = -60041154
$ws.Range("H70").Value = 0
$ws.Range("I70").Value = 0
$ws.Range("K70").Value = 0
$ws.Range("M70").ClearContents()
$ws.Range("H73").Value = 0
$ws.Range("I73").Value = 0
$ws.Range("K73").Value = 0
$ws.Range("M73").ClearContents()
$ws.Range("H92").Value = 1362.2
$ws.Range("I92").Value = 0
$ws.Range("J92").Value = 1362.2
$ws.Range("K92").Value = 0
$ws.Range("L92").Value = 4086.6
$ws.Range("M92").ClearContents()
$ws.Range("N92").Value = -6582.6
$ws.Range("H107").Value = 881.2308
$ws.Range("J107").Value = 935.6
$ws.Range("L107").Value = 2806.8
$ws.Range("N107").Value = -6646.8
$ws.Range("H113").Value = 2503.087
$ws.Range("J113").Value = 3224.3125
$ws.Range("L113").Value = 9672.9375
$ws.Range("N113").Value = -14012.9375

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 1333527.4
$ws.Range("I2").Value = 104.4
$ws.Range("J2").Value = 4000373.5
$ws.Range("K2").Value = 104.4
$ws.Range("L2").Value = 4000373.5
$ws.Range("M2").Value = 8.599999999999994
$ws.Range("N2").Value = -4000599.5
$ws.Range("H113").Value = 7170.7144
$ws.Range("I113").Value = 3496.6667
$ws.Range("J113").Value = 8172.727
$ws.Range("K113").Value = 3496.6667
$ws.Range("L113").Value = 8172.727
$ws.Range("M113").Value = -1326.6667
$ws.Range("N113").Value = -12512.727
$ws.Range("H132").Value = 10969.75
$ws.Range("I132").Value = 1940.75
$ws.Range("J132").Value = 19998.75
$ws.Range("K132").Value = 5822.25
$ws.Range("L132").Value = 59996.25
$ws.Range("M132").Value = -3292.25
$ws.Range("N132").Value = -65056.25
$ws.Range("H134").Value = 100000
$ws.Range("J134").Value = 100000
$ws.Range("L134").Value = 300000
$ws.Range("N134").Value = -305070

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5330.2
$ws.Range("I7").Value = 3993.625
$ws.Range("K7").Value = 3993.625
$ws.Range("M7").Value = -3881.625
$ws.Range("H22").Value = 4914.1177
$ws.Range("I22").Value = 2309.3333
$ws.Range("J22").Value = 9121.846
$ws.Range("K22").Value = 2309.3333
$ws.Range("L22").Value = 9121.846
$ws.Range("M22").Value = -2014.3333
$ws.Range("N22").Value = -9711.846
$ws.Range("H27").Value = 4914.1177
$ws.Range("I27").Value = 2309.3333
$ws.Range("J27").Value = 9121.846
$ws.Range("K27").Value = 2309.3333
$ws.Range("L27").Value = 9121.846
$ws.Range("M27").Value = -2202.3333
$ws.Range("N27").Value = -9335.846
$ws.Range("H55").Value = 406.42856
$ws.Range("I55").Value = 41.75
$ws.Range("J55").Value = 552.3
$ws.Range("K55").Value = 41.75
$ws.Range("L55").Value = 552.3
$ws.Range("M55").Value = 131.25
$ws.Range("N55").Value = -898.3
$ws.Range("H126").Value = 5330.2
$ws.Range("I126").Value = 3993.625
$ws.Range("K126").Value = 11980.875
$ws.Range("M126").Value = -9510.875
$ws.Range("H132").Value = 7280.184
$ws.Range("I132").Value = 4356.5
$ws.Range("K132").Value = 13069.5
$ws.Range("M132").Value = -10539.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 1524.6207
$ws.Range("I113").Value = 1251.75
$ws.Range("J113").Value = 1860.4615
$ws.Range("K113").Value = 3755.25
$ws.Range("L113").Value = 5581.3845
$ws.Range("M113").Value = -1585.25
$ws.Range("N113").Value = -9921.3845
$ws.Range("H132").Value = 88747.5
$ws.Range("I132").Value = 34996.668
$ws.Range("J132").Value = 250000
$ws.Range("K132").Value = 104990.004
$ws.Range("L132").Value = 750000
$ws.Range("M132").Value = -102460.004
$ws.Range("N132").Value = -755060
$ws.Range("H139").Value = 88639.875
$ws.Range("J139").Value = 88639.875
$ws.Range("L139").Value = 88639.875
$ws.Range("N139").Value = -98919.875
